$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "essence"
$ws.Range("B21").Value = "Esencia"
$ws.Range("C21").Value = "Essence"

$ws.Range("C22").Select()
